# "revisión sett recr dev con y sin flotas efecBo"
# Populate the new G/H "Label / model1" mini-table (rows 28-41) and the
# matching index column F (rows 29-41) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("G28").Value = "Label"
$ws.Range("H28").Value = "model1"

# Data rows: index (F), label (G), model1 value (H)
$rows = @(
  @{ r = 29; f = 1;  g = "TOTAL_like";             h = 156.92500000000001 },
  @{ r = 30; f = 2;  g = "Survey_like";             h = 16.3627 },
  @{ r = 31; f = 3;  g = "Age_comp_like";           h = 130.047 },
  @{ r = 32; f = 4;  g = "Parm_priors_like";        h = 0.97024299999999997 },
  @{ r = 33; f = 5;  g = "Recr_Virgin_billions";    h = 2.7384499999999998 },
  @{ r = 34; f = 6;  g = "SR_LN(R0)";               h = 14.822900000000001 },
  @{ r = 35; f = 7;  g = "SR_RkrPower_steep";       h = 0.38007200000000002 },
  @{ r = 36; f = 8;  g = "NatM_p_1_Fem_GP_1";       h = 0.34071400000000002 },
  @{ r = 37; f = 9;  g = "L_at_Amax_Fem_GP_1";      h = 65 },
  @{ r = 38; f = 10; g = "VonBert_K_Fem_GP_1";      h = 0.15 },
  @{ r = 39; f = 11; g = "SSB_Virgin_thousand_mt";  h = 2854.13 },
  @{ r = 40; f = 12; g = "Bratio_2017";             h = 0.36743100000000001 },
  @{ r = 41; f = 13; g = "SPRratio_2016";           h = 0.46838800000000003 }
)

foreach ($row in $rows) {
  $ws.Cells.Item($row.r, 6).Value = $row.f
  $ws.Cells.Item($row.r, 7).Value = $row.g
  $ws.Cells.Item($row.r, 8).Value = $row.h
}

# The new labels in column G are much wider than the old column, so widen it
# to fit the longest entry (mirrors the author resizing col G by hand).
$ws.Columns("G:G").AutoFit()

# Move the view / selection the way the author left it.
$ws.Range("H32").Select()
